$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 and Row 9 swap: Cardano <-> Dogecoin positions
$ws.Range("B8").Value = "Dogecoin"
$ws.Range("C8").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07677"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.24%  "

$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3065"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.54%  "

# Price / Volume updates
# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.439.34"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.79%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.854.22"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.17%  "

# Row 4
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.01%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6947"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.56%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.52"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.34%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07776"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.57%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.144"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.28%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.845.11"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.42%  "

# Row 14
$ws.Range("E14").Value = "  +0.57%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6917"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.72%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.276"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.44%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "29.429.12"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.70%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008337"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.61%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.101.02"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.11%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "237.64"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.21%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.72"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.08%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.0000"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.06%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.598"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.45%  "

# Row 24
$ws.Range("E24").Value = "  +0.02%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1489"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.64%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "159.78"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.94%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.875"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.83%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.23"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.22%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.527"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.03%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.234"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.63%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.145"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.01%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.208"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.24%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05096"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.13%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7726"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.02%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.882"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.43%  "

# Row 36
$ws.Range("E36").Value = "  +0.42%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.681"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.16%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.326.94"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +8.34%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01871"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.24%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.718"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.74%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9509"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.97%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "106.27"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.41%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.786"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.01%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.001"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.16%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "9.812"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.72%  "

# Row 46
$ws.Range("E46").Value = "  +2.68%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.999.19"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.98%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5219"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.89%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.784"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.00%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "62.99"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.81%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.952"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.76%  "
